# Add a new data row (2025-05-15) to the "Prices" sheet, right after the
# existing last row (74). All values in this sheet are stored as plain
# text (inline strings), including numbers like "37.5" or "5,343", so we
# force the new cells to text format before writing, then reset the
# number format back to the default "Normal" style to avoid leaving the
# cells with an explicit non-default style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 75
$rng = $ws.Range("A" + $row + ":J" + $row)
$rng.NumberFormat = "@"

$ws.Cells.Item($row, 1).Value  = "2025-05-15"
$ws.Cells.Item($row, 2).Value  = "37.5"
$ws.Cells.Item($row, 3).Value  = "37"
$ws.Cells.Item($row, 4).Value  = "0.98"
$ws.Cells.Item($row, 5).Value  = "0.265"
$ws.Cells.Item($row, 6).Value  = "0.09"
$ws.Cells.Item($row, 7).Value  = "5,343"
$ws.Cells.Item($row, 8).Value  = "7,998"
$ws.Cells.Item($row, 9).Value  = "8,048"
$ws.Cells.Item($row, 10).Value = "7.2236"

$rng.Style = "Normal"
